$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "English" "Английский"
Replace-Text " / Portuguese / French / Thai / Vietnamese / Spanish" " / Португальский / Французский / Тайский / Вьетнамский / Испанский"
Replace-Text "Brief" "Кратко"
Replace-Text "An email sent to partners in the target country who have RSVPed no. It will be sent via customer.io" "Письмо, отправленное партнерам в целевой стране, которые ответили `"нет`". Оно будет отправлено через customer.io"
Replace-Text "Target audience" "Целевая аудитория"
Replace-Text "Invited partners who RSVP no" "Приглашенные партнеры, которые ответили `"нет`""
Replace-Text "Subject line" "Тема письма"
Replace-Text ": Thinking of you at " ": Думаем о вас на "
Replace-Text "[EVENT NAME]" "[НАЗВАНИЕ МЕРОПРИЯТИЯ]"
Replace-Text "We’ll miss you at the " "Мы будем скучать по вам на "
Replace-Text "Dear " "Здравствуйте, "
Replace-Text "[PARTNER NAME]" "[ИМЯ ПАРТНЕРА]"
Replace-Text "Thank you for taking the time to respond to our invitation to the upcoming " "Спасибо, что нашли время ответить на наше приглашение на грядущее "
Replace-Text ". We were really looking forward to seeing you there." ". Мы действительно с нетерпением ждали встречи с вами."
Replace-Text "Even though we’re disappointed we can’t meet you, we understand that scheduling conflicts and other commitments sometimes come up. " "Хотя мы разочарованы, что не можем встретиться с вами, мы понимаем, что иногда возникают конфликты в расписании и другие обязательства."
Replace-Text "If you’re comfortable sharing it with us, we’d like to know why you responded no. Please reply to this email as your feedback could help us make improvements in our event planning processes and better serve you in the future." "Если вы не против поделиться с нами, нам хотелось бы узнать, почему вы ответили `"нет`". Пожалуйста, ответьте на это письмо, так как ваш отзыв может помочь нам улучшить процессы планирования наших мероприятий и лучше служить вам в будущем."
Replace-Text "We hope to see you at our future events. " "Мы надеемся увидеть вас на наших будущих мероприятиях."
Replace-Text "If you have any questions, please contact us via " "Если у вас есть вопросы, пожалуйста, свяжитесь с нами через "
Replace-Text "live chat" "чат"
Replace-Text " or " " или "
Replace-Text "If you have any questions, please contact your country manager, " "Если у вас есть вопросы, пожалуйста, свяжитесь с вашим региональным менеджером, "
Replace-Text "[NAME]" "[ИМЯ]"
Replace-Text ", at " ", по адресу "
Replace-Text "[EMAIL ADDRESS]" "[АДРЕС ЭЛЕКТРОННОЙ ПОЧТЫ]"
Replace-Text "[WHATSAPP NO]" "[НОМЕР WHATSAPP]"
Replace-Text "choose either one" "выберите один из них"
